$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3 ("Shivam") becomes "Sachin" and a matching value is filled into the
# previously-empty A3 cell (no rows are inserted/shifted - rows 4-7 stay put).
$ws.Range("A3").Value = "Sachin"
$ws.Range("B3").Value = "Sachin"
